$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product_name column (B) from parka1/parka2/parka3 to pants1/pants2/pants3
$ws.Range("B2").Value = "pants1"
$ws.Range("B3").Value = "pants2"
$ws.Range("B4").Value = "pants3"

# Update size column (M) from S to L
$ws.Range("M2").Value = "L"
$ws.Range("M3").Value = "L"
$ws.Range("M4").Value = "L"

# Update the active cell selection
$ws.Range("K8").Select()
